$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H10").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 0
